$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New rows (276-281): six new "St.Petersburg, Russia" live-camera entries
# appended to the bottom of the list. Columns: A=Category, B=lat,long,
# C=Title, D=City, E=Country, F=YouTube video id, G=Channel (hyperlinked).
# ---------------------------------------------------------------------------

$rowsData = @(
    @{ Row=276; A="LIVE, RIVER, WARSHIP";          B="59.9552456170424, 30.337814739753615";   C="LIVE CAMERA Russian Cruiser Aurora St.Petersburg Russia Россия онлайн Санкт-Петербург крейсер Аврора"; D="St.Petersburg"; E="Russia"; F="JYLJDYxV_ms";  G="Mobotix Webcams Russia Intro видеонаблюдение в Санкт-Петербурге" },
    @{ Row=277; A="LIVE, TRAFFIC, STREET";         B="59.934273683681106, 30.335150370158708"; C="LIVE Nevskiy avenue St. Petersburg Russia, Gostiny Dvor. Невский пр. Санкт-Петербург, Гостиный двор"; D="St.Petersburg"; E="Russia"; F="h1wly909BYw";  G="Mobotix Webcams Russia Intro видеонаблюдение в Санкт-Петербурге" },
    @{ Row=278; A="LIVE, TRAFFIC, STREET";         B="60.00297893433313, 30.297868997414106";  C="🕹️ PTZ LIVE CAMERA St. Petersburg, Pionerskaya. Метро Пионерская пр. Испытателей и Коломяжский пр."; D="St.Petersburg"; E="Russia"; F="hbPntLL5eU8";  G="Mobotix Webcams Russia Intro видеонаблюдение в Санкт-Петербурге" },
    @{ Row=279; A="LIVE, TRAFFIC, STREET";         B="60.000633038756554, 30.255864844735935"; C="CROSSROADS Camera Gakkelevskaya st. Bogatyrsky ave., St. Petersburg Онлайн камера на перекрестке СПб"; D="St.Petersburg"; E="Russia"; F="ZlDohRExM-A";  G="Mobotix Webcams Russia Intro видеонаблюдение в Санкт-Петербурге" },
    @{ Row=280; A="LIVE, RIVER, BRIDGE, TRAFFIC";  B="59.94004991910285, 30.303912259408182";  C="Дворцовая набережная и река Нева в прямом эфире. Palace embankment and Neva river ship cam online"; D="St.Petersburg"; E="Russia"; F="6qjTpNw-psE";  G="Mobotix Webcams Russia Intro видеонаблюдение в Санкт-Петербурге" },
    @{ Row=281; A="LIVE, TRAFFIC, STREET";         B="59.92919635814432, 30.34401868413709";   C="4K video LIVE CAMERA Saint Petersburg, Russia. Rubinstein Street Улица Рубинштейна онлайн камера"; D="St.Petersburg"; E="Russia"; F="GIUTYf0Fpic";  G="Mobotix Webcams Russia Intro видеонаблюдение в Санкт-Петербурге" }
)

foreach ($r in $rowsData) {
    $row = $r.Row
    $ws.Range("A$row").Value = $r.A
    $ws.Range("B$row").Value = $r.B
    $ws.Range("C$row").Value = $r.C
    $ws.Range("D$row").Value = $r.D
    $ws.Range("E$row").Value = $r.E
    $ws.Range("F$row").Value = $r.F
    $ws.Range("G$row").Value = $r.G
}

# ---------------------------------------------------------------------------
# Formatting: match the look of the rows directly above (257-275).
#   - columns A & E use the thin-border "row" style
#   - column G carries the hyperlink (and picks up the Hyperlink look)
# Apply the hyperlinks FIRST, then copy the existing bordered-row formatting
# over A/E/G so every new row matches the established table styling.
# ---------------------------------------------------------------------------

foreach ($r in $rowsData) {
    $row = $r.Row
    $ws.Hyperlinks.Add($ws.Range("G$row"), "https://www.youtube.com/@msbud2", "", "", $r.G) | Out-Null
}

$ws.Range("A275").Copy() | Out-Null
foreach ($r in $rowsData) {
    $row = $r.Row
    $ws.Range("A$row").PasteSpecial(-4122) | Out-Null
}

$ws.Range("E275").Copy() | Out-Null
foreach ($r in $rowsData) {
    $row = $r.Row
    $ws.Range("E$row").PasteSpecial(-4122) | Out-Null
}

$ws.Range("G275").Copy() | Out-Null
foreach ($r in $rowsData) {
    $row = $r.Row
    $ws.Range("G$row").PasteSpecial(-4122) | Out-Null
}

$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Tidy up the existing rows above (257-275): re-apply their own formatting so
# the workbook's style table collapses the redundant "bordered" style onto a
# single canonical entry, same as Excel does when it re-saves the sheet.
# ---------------------------------------------------------------------------

$ws.Range("A257").Copy() | Out-Null
$ws.Range("A257:A275").PasteSpecial(-4122) | Out-Null
$ws.Range("E257").Copy() | Out-Null
$ws.Range("E257:E275").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# Update the view: the previous selection/frozen-pane anchor tracked the old
# last row (A276); move it down to the new last row (A284).
# ---------------------------------------------------------------------------

$ws.Range("A284").Select() | Out-Null

$wb.Save()
